{"js": "const replacements = [\n  [\"668\u00d76=4008\", \"383\u00d74=1532\"],\n  [\"910\u00d72=1820\", \"560\u00d73=1680\"],\n  [\"227\u00d74=908\", \"341\u00d75=1705\"],\n  [\"242\u00d78=1936\", \"396\u00d74=1584\"],\n  [\"232\u00d74=928\", \"205\u00d73=615\"],\n  [\"217\u00d74=868\", \"801\u00d79=7209\"],\n  [\"993\u00d76=5958\", \"301\u00d75=1505\"],\n  [\"468\u00d75=2340\", \"208\u00d73=624\"],\n  [\"751\u00d75=3755\", \"999\u00d79=8991\"],\n  [\"664\u00d72=1328\", \"651\u00d78=5208\"],\n  [\"495\u00d72=990\", \"937\u00d77=6559\"],\n  [\"769\u00d78=6152\", \"113\u00d77=791\"],\n  [\"675\u00d76=4050\", \"622\u00d73=1866\"],\n  [\"640\u00d75=3200\", \"336\u00d79=3024\"],\n  [\"994\u00d78=7952\", \"558\u00d79=5022\"],\n  [\"161\u00d74=644\", \"601\u00d76=3606\"],\n  [\"669\u00d73=2007\", \"574\u00d79=5166\"],\n  [\"433\u00d79=3897\", \"356\u00d78=2848\"],\n  [\"224\u00d79=2016\", \"808\u00d72=1616\"],\n  [\"295\u00d79=2655\", \"803\u00d75=4015\"],\n  [\"708\u00d76=4248\", \"885\u00d72=1770\"],\n  [\"972\u00d79=8748\", \"856\u00d78=6848\"],\n  [\"207\u00d77=1449\", \"948\u00d78=7584\"],\n  [\"826\u00d75=4130\", \"164\u00d73=492\"],\n  [\"636\u00d76=3816\", \"648\u00d75=3240\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"668\u00d76=4008\", \"383\u00d74=1532\"),\n    @(\"910\u00d72=1820\", \"560\u00d73=1680\"),\n    @(\"227\u00d74=908\", \"341\u00d75=1705\"),\n    @(\"242\u00d78=1936\", \"396\u00d74=1584\"),\n    @(\"232\u00d74=928\", \"205\u00d73=615\"),\n    @(\"217\u00d74=868\", \"801\u00d79=7209\"),\n    @(\"993\u00d76=5958\", \"301\u00d75=1505\"),\n    @(\"468\u00d75=2340\", \"208\u00d73=624\"),\n    @(\"751\u00d75=3755\", \"999\u00d79=8991\"),\n    @(\"664\u00d72=1328\", \"651\u00d78=5208\"),\n    @(\"495\u00d72=990\", \"937\u00d77=6559\"),\n    @(\"769\u00d78=6152\", \"113\u00d77=791\"),\n    @(\"675\u00d76=4050\", \"622\u00d73=1866\"),\n    @(\"640\u00d75=3200\", \"336\u00d79=3024\"),\n    @(\"994\u00d78=7952\", \"558\u00d79=5022\"),\n    @(\"161\u00d74=644\", \"601\u00d76=3606\"),\n    @(\"669\u00d73=2007\", \"574\u00d79=5166\"),\n    @(\"433\u00d79=3897\", \"356\u00d78=2848\"),\n    @(\"224\u00d79=2016\", \"808\u00d72=1616\"),\n    @(\"295\u00d79=2655\", \"803\u00d75=4015\"),\n    @(\"708\u00d76=4248\", \"885\u00d72=1770\"),\n    @(\"972\u00d79=8748\", \"856\u00d78=6848\"),\n    @(\"207\u00d77=1449\", \"948\u00d78=7584\"),\n    @(\"826\u00d75=4130\", \"164\u00d73=492\"),\n    @(\"636\u00d76=3816\", \"648\u00d75=3240\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}\n"}
